$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Base the new row 47 on row 44, which already carries the "reviewed / wrong-mark"
# formatting pattern used for finished entries (A/B/C/D/E/F/H columns).
$ws.Range("A44:H44").Copy($ws.Range("A47:H47"))

# G column should look like a "not yet reviewed" entry (same as G45/G46) rather
# than the reviewed-date style copied from G44.
$ws.Range("G46").Copy($ws.Range("G47"))

# --- Fill in the new problem's data -----------------------------------------
$ws.Range("A47").Value2 = "376. Wiggle Subsequence"
$ws.Range("C47").Value2 = "https://leetcode.com/problems/wiggle-subsequence/"
$ws.Range("D47").Value2 = 44551
$ws.Range("E47").Value2 = "动态规划，摇摆序列"

$f47 = $ws.Range("F47")
$f47.Value2 = "纯dp，数组要根据末尾两元素单调性分类；可以用贪心优化到O(n)"
$f47.Font.Name = "Times New Roman"
$f47.Characters(1,1).Font.Name = "宋体"
$f47.Characters(2,2).Font.Name = "Times New Roman"
$f47.Characters(4,25).Font.Name = "宋体"
$f47.Characters(29,4).Font.Name = "Times New Roman"

# Row height matches the other two-line entries in this part of the table.
$ws.Rows.Item(47).RowHeight = 42

# The mistake marker in H47 uses a non-bold variant of the usual font.
$ws.Range("H47").Font.Bold = $false

# Link the new problem's title cell to LeetCode.
$ws.Hyperlinks.Add($ws.Range("C47"), "https://leetcode.com/problems/wiggle-subsequence/")

# Reflect the cursor/viewport position from the saved workbook.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 34
$ws.Range("I47").Select() | Out-Null
